$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.124546527862549
$ws.Range("B1").Value = 2.272533416748047
$ws.Range("C1").Value = 10.25321006774902
$ws.Range("D1").Value = 1.943143606185913
$ws.Range("E1").Value = 1.285963654518127
